$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.138.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.60%  '

$ws.Range("D3").Value = "'2.907.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'566.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.36%  '

$ws.Range("D6").Value = "'143.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.07%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = "'2.905.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -1.09%  '

$ws.Range("D10").Value = "'6.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.67%  '

$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("D14").Value = "'32.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.15%  '

$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("D16").Value = "'3.390.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("D17").Value = "'62.070.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = "'2.920.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = "'6.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").Value = "'430.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").Value = "'13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.00%  '

$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("E23").Value = '  -0.96%  '

$ws.Range("D24").Value = "'78.61"
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  +1.15%  '

$ws.Range("D26").Value = "'10.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  -2.65%  '

$ws.Range("E29").Value = '  +5.19%  '

$ws.Range("E30").Value = '  -3.65%  '

$ws.Range("E31").Value = '  -2.10%  '

$ws.Range("E32").Value = '  -4.09%  '

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.14%  '

$ws.Range("D34").Value = "'25.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.89%  '

$ws.Range("E35").Value = '  -3.28%  '

$ws.Range("D36").Value = "'0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.46%  '

$ws.Range("E37").Value = '  -2.25%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = "'2.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = "'48.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("E40").Value = '  -4.92%  '

$ws.Range("E41").Value = '  -0.76%  '

$ws.Range("E42").Value = '  -2.50%  '

$ws.Range("D43").Value = "'40.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.87%  '

$ws.Range("D44").Value = "'2.717.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("E45").Value = '  -2.26%  '

$ws.Range("D46").Value = "'133.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").Value = "'345.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").Value = "'0.000221"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.54%  '

$ws.Range("E51").Value = '  -0.77%  '
